$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 3.8
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4.5
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("T3").Value = 2.37
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AE3").Value = 19
$ws.Range("AG3").Value = 8.5
$ws.Range("AI3").Value = 15
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AX3").Value = 23
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 351

# Row 4
$ws.Range("S4").Value = 1.54

# Row 5
$ws.Range("I5").Value = 3.6
$ws.Range("N5").Value = 8
$ws.Range("S5").Value = 1.47
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67

# Row 6
$ws.Range("S6").Value = 1.47

# Row 9
$ws.Range("Q9").Value = 1.95
$ws.Range("R9").Value = 1.9
